# Regenerate the localization-status report: the handoff that was sitting
# "Ready for handoff" has moved on to "In Translation". Update the status
# cell wherever it's reported (Overview roll-up + each locale sheet) and
# let the column shrink to fit the shorter text, just like Excel does when
# the report is regenerated.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Overview sheet: the zh-cn/de-de status columns (E, F) were sized to fit
# "Ready for handoff" and now fit the shorter "In Translation" text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.43
$wsOverview.Columns.Item(6).ColumnWidth = 12.43

# Per-locale sheets: the Status column (C) shrinks the same way.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.43

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.43
